$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "清洁度"
$ws.Range("A3").Value = "霉菌"
$ws.Range("A4").Value = "滴虫"
$ws.Range("A5").Value = "线索细胞！"
$ws.Range("A6").Value = "AST/ALT"
$ws.Range("E6").Value = "-"
$ws.Range("A7").Value = "鳞状上皮（鳞状上皮细胞"
